# Saldo_guide.xlsx update:
#  - Bump the reference date (column G, "Dt. Referencia") for every data row
#    from 45434 (2024-05-22) to 45435 (2024-05-23).
#  - Update the balance figures for account row 117 (D and H columns,
#    "Saldo Previsto" / "Vl. Total") from 170.16 to 22358.59.
#  - Rename the worksheet to reflect the new export file name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Dt. Referencia" column for all data rows (2 through 257).
$ws.Range("G2:G257").Value = 45435

# Update the balance values on row 117.
$ws.Cells.Item(117, 4).Value = 22358.59
$ws.Cells.Item(117, 8).Value = 22358.59

# Rename the sheet to match the new export name.
$ws.Name = "IClientBalance-20240523-094816-"
